$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B25").Value = 6482
$ws.Range("D25").Value = 6044863
$ws.Range("E25").Value = 932.5614008022216
$ws.Range("F25").Value = 10.03225258869462
$ws.Range("H25").Value = 26.59807248233765
